# Update Daily Report: adds the new business-day snapshot (2026-01-16, serial 46038)
# to Daily_Data, and refreshes the derived Today_Summary / Monthly_Stats sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Daily_Data: append 22 new rows (rows 222-243) for date serial 46038
# ---------------------------------------------------------------------------
$daily = $wb.Worksheets.Item("Daily_Data")

$newDateSerial = 46038
$dateNumberFormat = "YYYY-MM-DD HH:MM:SS"

# Each entry: Label, PREV_TOTAL(C), RECEIVED(D), WITHDRAWN(E), NET_CHANGE(F), ADJUSTMENT(G), TOTAL_TODAY(H)
$newRows = @(
    @("ASAHI DEPOSITORY LLC Registered", 0, 0, 0, 0, 0, 0),
    @("ASAHI DEPOSITORY LLC Eligible", 0, 0, 0, 0, 0, 0),
    @("BRINK'S, INC. Registered", 95517.77499999999, 0, 0, 0, -3784.014, 91733.761),
    @("BRINK'S, INC. Eligible", 23710.274, 0, 0, 0, 3784.014, 27494.288),
    @("CNT DEPOSITORY, INC. Registered", 1246.06, 0, 0, 0, 0, 1246.06),
    @("CNT DEPOSITORY, INC. Eligible", 0, 0, 0, 0, 0, 0),
    @("DELAWARE DEPOSITORY Registered", 1633.941, 0, 0, 0, 0, 1633.941),
    @("DELAWARE DEPOSITORY Eligible", 18459.584, 0, 0, 0, 0, 18459.584),
    @("HSBC BANK, USA Registered", 1295.223, 0, 0, 0, 0, 1295.223),
    @("HSBC BANK, USA Eligible", 9381.513999999999, 0, 0, 0, 0, 9381.513999999999),
    @("INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered", 2395.448, 0, 0, 0, 0, 2395.448),
    @("INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible", 0, 0, 0, 0, 0, 0),
    @("JP MORGAN CHASE BANK NA Registered", 124991.729, 0, 0, 0, -10006.15, 114985.579),
    @("JP MORGAN CHASE BANK NA Eligible", 125407.673, 0, 0, 0, 10006.15, 135413.823),
    @("LOOMIS INTERNATIONAL (US) LLC Registered", 68084.33, 0, 0, 0, -15711.682, 52372.648),
    @("LOOMIS INTERNATIONAL (US) LLC Eligible", 116365.524, 0, 0, 0, 15711.682, 132077.206),
    @("MALCA-AMIT USA, LLC Registered", 395.145, 0, 0, 0, 0, 395.145),
    @("MALCA-AMIT USA, LLC Eligible", 0, 0, 0, 0, 0, 0),
    @("MANFRA, TORDELLA & BROOKES, LLC Registered", 60301.249, 0, 0, 0, -10080.829, 50220.42),
    @("MANFRA, TORDELLA & BROOKES, LLC Eligible", 1068.408, 0, 0, 0, 10080.829, 11149.237),
    @("STONEX PRECIOUS METALS LLC Registered", 14122.765, 0, 0, 0, 0, 14122.765),
    @("STONEX PRECIOUS METALS LLC Eligible", 16.075, 0, 0, 0, 0, 16.075)
)

$startRow = 222
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $daily.Cells.Item($r, 1).Value = $newDateSerial
    $daily.Cells.Item($r, 1).NumberFormat = $dateNumberFormat
    $daily.Cells.Item($r, 2).Value = $data[0]
    $daily.Cells.Item($r, 3).Value = $data[1]
    $daily.Cells.Item($r, 4).Value = $data[2]
    $daily.Cells.Item($r, 5).Value = $data[3]
    $daily.Cells.Item($r, 6).Value = $data[4]
    $daily.Cells.Item($r, 7).Value = $data[5]
    $daily.Cells.Item($r, 8).Value = $data[6]
}

# ---------------------------------------------------------------------------
# 2. Today_Summary: refresh Eligible/Registered/Total_Stock for companies
#    whose latest-day figures changed (B=Eligible, C=Registered, D=Total)
# ---------------------------------------------------------------------------
$today = $wb.Worksheets.Item("Today_Summary")

$todayUpdates = @{
    3  = @(27494.288, 91733.761, 119228.049)    # BRINK'S, INC.
    8  = @(135413.823, 114985.579, 250399.402)  # JP MORGAN CHASE BANK NA
    9  = @(132077.206, 52372.648, 184449.854)   # LOOMIS INTERNATIONAL (US) LLC
    11 = @(11149.237, 50220.42, 61369.657)      # MANFRA, TORDELLA & BROOKES, LLC
}

foreach ($r in $todayUpdates.Keys) {
    $vals = $todayUpdates[$r]
    $today.Cells.Item($r, 2).Value = $vals[0]
    $today.Cells.Item($r, 3).Value = $vals[1]
    $today.Cells.Item($r, 4).Value = $vals[2]
}

# ---------------------------------------------------------------------------
# 3. Monthly_Stats: refresh month-to-date rollups
# ---------------------------------------------------------------------------
$monthly = $wb.Worksheets.Item("Monthly_Stats")

# Top summary row (row 2): Eligible / Registered / Grand_Total
$monthly.Cells.Item(2, 2).Value = 333991.727
$monthly.Cells.Item(2, 3).Value = 330400.99
$monthly.Cells.Item(2, 4).Value = 664392.7169999999

# Detail rows (E column = TOTAL_TODAY, i.e. latest snapshot for the company/type)
$monthlyDetailUpdates = @{
    9  = 27494.288    # BRINK'S, INC. Eligible
    10 = 91733.761    # BRINK'S, INC. Registered
    19 = 135413.823   # JP MORGAN CHASE BANK NA Eligible
    20 = 114985.579   # JP MORGAN CHASE BANK NA Registered
    21 = 132077.206   # LOOMIS INTERNATIONAL (US) LLC Eligible
    22 = 52372.648    # LOOMIS INTERNATIONAL (US) LLC Registered
    25 = 11149.237    # MANFRA, TORDELLA & BROOKES, LLC Eligible
    26 = 50220.42     # MANFRA, TORDELLA & BROOKES, LLC Registered
}

foreach ($r in $monthlyDetailUpdates.Keys) {
    $monthly.Cells.Item($r, 5).Value = $monthlyDetailUpdates[$r]
}

Write-Host "Applied Daily Report update for 2026-01-16 (serial 46038)"
